# Update oyster spreadsheet for sheet names:
#   "oyster height"    -> "excavation length"
#   "sample frequency" -> "excavation count"
#
# This affects:
#   1) The two worksheet (tab) names themselves.
#   2) The "sheet name" column (column F) on the "glossary" worksheet,
#      which lists, per glossary field, the sheets that use that field
#      (comma-separated lists of sheet names).

$wb = $excel.ActiveWorkbook

# --- 1) Rename the worksheets -------------------------------------------
$wsOysterHeight = $wb.Worksheets.Item("oyster height")
$wsOysterHeight.Name = "excavation length"

$wsSampleFrequency = $wb.Worksheets.Item("sample frequency")
$wsSampleFrequency.Name = "excavation count"

# --- 2) Update references on the glossary sheet --------------------------
$glossary = $wb.Worksheets.Item("glossary")

$usedRange = $glossary.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $glossary.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string]) {
        if ($val.Contains("oyster height") -or $val.Contains("sample frequency")) {
            $newVal = $val.Replace("oyster height", "excavation length")
            $newVal = $newVal.Replace("sample frequency", "excavation count")
            $cell.Value = $newVal
        }
    }
}
